# [monitoring_dift_drop] Change from monitor to monitoring core.
#
# The "assembly" diagram on slide 3 has a small textbox ("TextBox 8",
# shape id 9) whose whole text is just "Monitor". Rename it to
# "Monitoring Core".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

$target = $null

# Prefer an exact-text match so this keeps working even if shape
# z-order/indices ever shift around.
for ($j = 1; $j -le $s.Shapes.Count; $j++) {
    $shp = $s.Shapes.Item($j)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.HasText -and $shp.TextFrame.TextRange.Text -eq "Monitor") {
            $target = $shp
            break
        }
    }
}

# Fall back to the known shape name, then to the known index.
if ($target -eq $null) {
    try {
        $target = $s.Shapes.Item("TextBox 8")
    } catch {
        $target = $s.Shapes.Item(6)
    }
}

$target.TextFrame.TextRange.Text = "Monitoring Core"
